$wb = $excel.ActiveWorkbook

$wsResults = $wb.Worksheets.Item("Results")
$wsInstructions = $wb.Worksheets.Item("Instructions")

# Results!O3 : "Grab" -> "Grab-MassWateR"
$wsResults.Range("O3").Value = "Grab-MassWateR"

# Instructions!C20 : "Grab" -> "Grab-MassWateR"
$wsInstructions.Range("C20").Value = "Grab-MassWateR"

# Instructions!B20 : update instructional text to reference "Grab-MassWateR" (and siblings)
$wsInstructions.Range("B20").Value = "For WQX:  Enter the method ID used for this sample collection.  Not applicable for field measurement/observations.  Method IDs are defined in WQX by organization.  MassWateR will assign a default value of ""Grab-MassWateR"" if nothing is entered, but this requires a Method Context of ""MassWateR"" in the WQXMeta file.  Standard method IDs that can be used by any organization under the MassWateR context are ""Grab-MassWateR"", ""Pole-MassWateR"", and ""Basket-MassWateR""."

# Restore / set the final active-cell selections to match the saved view state
$wsInstructions.Activate()
$wsInstructions.Range("A6").Select()

$wsResults.Activate()
$wsResults.Range("A3").Select()
